# Scheduled runner update: refresh market-price derived columns
# (currentAveragePrice / NQ / HQ / LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ)
# across the ALC, ARM, BSM, CRP, CUL, GSM and LTW sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 35714564
$ws.Range("I33").Value = 35714564
$ws.Range("K33").Value = 35714564
$ws.Range("M33").Value = -35714335
$ws.Range("H69").Value = 4385
$ws.Range("I69").Value = 5125
$ws.Range("J69").Value = 3398.3333
$ws.Range("K69").Value = 15375
$ws.Range("L69").Value = 10194.9999
$ws.Range("M69").Value = -14501
$ws.Range("N69").Value = -11942.9999
$ws.Range("H72").Value = 4385
$ws.Range("I72").Value = 5125
$ws.Range("J72").Value = 3398.3333
$ws.Range("K72").Value = 46125
$ws.Range("L72").Value = 30584.9997
$ws.Range("M72").Value = -41757
$ws.Range("N72").Value = -39320.9997
$ws.Range("H116").Value = 4952.0625
$ws.Range("I116").Value = 6664.7144
$ws.Range("J116").Value = 3620
$ws.Range("K116").Value = 6664.7144
$ws.Range("L116").Value = 3620
$ws.Range("M116").Value = -3222.7144
$ws.Range("N116").Value = -10504

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1128.762
$ws.Range("I2").Value = 1013.6
$ws.Range("J2").Value = 1416.6666
$ws.Range("K2").Value = 1013.6
$ws.Range("L2").Value = 1416.6666
$ws.Range("M2").Value = -900.6
$ws.Range("N2").Value = -1642.6666
$ws.Range("H61").Value = 3377.2856
$ws.Range("I61").Value = 3230.7646
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 3230.7646
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -3018.7646
$ws.Range("N61").Value = -4424
$ws.Range("H63").Value = 41668772
$ws.Range("I63").Value = 50001930
$ws.Range("J63").Value = 3000
$ws.Range("K63").Value = 50001930
$ws.Range("L63").Value = 3000
$ws.Range("M63").Value = -50001244
$ws.Range("N63").Value = -4372
$ws.Range("H66").Value = 41668772
$ws.Range("I66").Value = 50001930
$ws.Range("J66").Value = 3000
$ws.Range("K66").Value = 250009650
$ws.Range("L66").Value = 15000
$ws.Range("M66").Value = -250006218
$ws.Range("N66").Value = -21864
$ws.Range("H74").Value = 2012.3125
$ws.Range("I74").Value = 1168.1305
$ws.Range("J74").Value = 4169.6665
$ws.Range("K74").Value = 1168.1305
$ws.Range("L74").Value = 4169.6665
$ws.Range("M74").Value = -294.1305
$ws.Range("N74").Value = -5917.6665
$ws.Range("H77").Value = 2012.3125
$ws.Range("I77").Value = 1168.1305
$ws.Range("J77").Value = 4169.6665
$ws.Range("K77").Value = 5840.6525
$ws.Range("L77").Value = 20848.3325
$ws.Range("M77").Value = -1472.6525
$ws.Range("N77").Value = -29584.3325
$ws.Range("H116").Value = 1128.762
$ws.Range("I116").Value = 1013.6
$ws.Range("J116").Value = 1416.6666
$ws.Range("K116").Value = 1013.6
$ws.Range("L116").Value = 1416.6666
$ws.Range("M116").Value = 1280.4
$ws.Range("N116").Value = -6004.6666
$ws.Range("H136").Value = 3377.2856
$ws.Range("I136").Value = 3230.7646
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 9692.293799999999
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -7142.293799999999
$ws.Range("N136").Value = -17100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1128.762
$ws.Range("I3").Value = 1013.6
$ws.Range("J3").Value = 1416.6666
$ws.Range("K3").Value = 1013.6
$ws.Range("L3").Value = 1416.6666
$ws.Range("M3").Value = -899.6
$ws.Range("N3").Value = -1644.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1309.8572
$ws.Range("I99").Value = 1232.4
$ws.Range("J99").Value = 1503.5
$ws.Range("K99").Value = 1232.4
$ws.Range("L99").Value = 1503.5
$ws.Range("M99").Value = 265.5999999999999
$ws.Range("N99").Value = -4499.5
$ws.Range("H126").Value = 1309.8572
$ws.Range("I126").Value = 1232.4
$ws.Range("J126").Value = 1503.5
$ws.Range("K126").Value = 3697.2
$ws.Range("L126").Value = 4510.5
$ws.Range("M126").Value = -1227.2
$ws.Range("N126").Value = -9450.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 4482.185
$ws.Range("I136").Value = 2274.4546
$ws.Range("J136").Value = 6000
$ws.Range("K136").Value = 6823.3638
$ws.Range("L136").Value = 18000
$ws.Range("M136").Value = -1723.3638
$ws.Range("N136").Value = -28200

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 18044.5
$ws.Range("J15").Value = 18044.5
$ws.Range("L15").Value = 18044.5
$ws.Range("N15").Value = -18620.5
$ws.Range("H80").Value = 15154293
$ws.Range("I80").Value = 33335994
$ws.Range("J80").Value = 2875
$ws.Range("K80").Value = 33335994
$ws.Range("L80").Value = 2875
$ws.Range("M80").Value = -33334996
$ws.Range("N80").Value = -4871
$ws.Range("H81").Value = 18044.5
$ws.Range("J81").Value = 18044.5
$ws.Range("L81").Value = 18044.5
$ws.Range("N81").Value = -20040.5
$ws.Range("H83").Value = 15154293
$ws.Range("I83").Value = 33335994
$ws.Range("J83").Value = 2875
$ws.Range("K83").Value = 166679970
$ws.Range("L83").Value = 14375
$ws.Range("M83").Value = -166674978
$ws.Range("N83").Value = -24359
$ws.Range("H84").Value = 18044.5
$ws.Range("J84").Value = 18044.5
$ws.Range("L84").Value = 54133.5
$ws.Range("N84").Value = -64117.5
$ws.Range("H97").Value = 1649.3636
$ws.Range("I97").Value = 1158.5714
$ws.Range("J97").Value = 2508.25
$ws.Range("K97").Value = 1158.5714
$ws.Range("L97").Value = 2508.25
$ws.Range("M97").Value = -662.5714
$ws.Range("N97").Value = -3500.25
$ws.Range("H126").Value = 2871.976
$ws.Range("I126").Value = 1772.3214
$ws.Range("J126").Value = 5071.2856
$ws.Range("K126").Value = 5316.9642
$ws.Range("L126").Value = 15213.8568
$ws.Range("M126").Value = -2846.9642
$ws.Range("N126").Value = -20153.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2116.1538
$ws.Range("I82").Value = 1968.3334
$ws.Range("J82").Value = 2242.8572
$ws.Range("K82").Value = 1968.3334
$ws.Range("L82").Value = 2242.8572
$ws.Range("M82").Value = -1607.3334
$ws.Range("N82").Value = -2964.8572
$ws.Range("H85").Value = 2116.1538
$ws.Range("I85").Value = 1968.3334
$ws.Range("J85").Value = 2242.8572
$ws.Range("K85").Value = 1968.3334
$ws.Range("L85").Value = 2242.8572
$ws.Range("M85").Value = -720.3334
$ws.Range("N85").Value = -4738.8572
